$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "33.839.28"
$ws.Range("E2").Value = "  +7.06%  "
$ws.Range("D3").Value = "1.777.91"
$ws.Range("E3").Value = "  +3.91%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.13"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.43%  "
$ws.Range("E6").Value = "  +4.36%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.24%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "30.63"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +2.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.58"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +2.72%  "
$ws.Range("E10").Value = "  +2.86%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0665"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +3.31%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0922"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.26%  "
$ws.Range("D13").Value = "2.034.14"
$ws.Range("E13").Value = "  +4.02%  "
$ws.Range("D14").Value = "1.776.76"
$ws.Range("E14").Value = "  +3.86%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.624"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.78%  "
$ws.Range("D16").Value = "33.799.70"
$ws.Range("E16").Value = "  +7.17%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "9.99"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -2.67%  "
$ws.Range("E18").Value = "  -0.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "68.53"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.86%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "251.16"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.17%  "
$ws.Range("E21").Value = "  +1.90%  "
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.27"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.66%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.17"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.53%  "
$ws.Range("E25").Value = "  -0.23%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.23"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.70%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.43"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.26%  "
$ws.Range("E28").Value = "  +0.84%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.97"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +2.27%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.26%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.80"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.64%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0514"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.98%  "
$ws.Range("E33").Value = "  +3.21%  "
$ws.Range("E34").Value = "  +4.57%  "
$ws.Range("E35").Value = "  +6.08%  "
$ws.Range("D36").Value = "1.477.37"
$ws.Range("E36").Value = "  -2.66%  "
$ws.Range("E37").Value = "  +2.57%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.630"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +3.29%  "
$ws.Range("E39").Value = "  +2.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "83.04"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.67%  "
$ws.Range("E41").Value = "  +1.75%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.69"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.885"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +3.32%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.08"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.52%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0507"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.46%  "
$ws.Range("E46").Value = "  +3.41%  "
$ws.Range("D47").Value = "1.929.72"
$ws.Range("E47").Value = "  +4.57%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.69"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.20%  "
$ws.Range("E49").Value = "  +0.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.85"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +15.37%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.00"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.74%  "
